# Generate Report for Handoff
#
# This script applies the "Ready for handoff" status/report-generation
# update to the localization-status workbook:
#   - Status cells move from "In Translation" to "Ready for handoff"
#   - The "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#     timestamps are refreshed
#   - The Status/date columns are widened slightly on every sheet

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" -------------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Refreshed handoff timestamps --------------------------------------
$overview.Range("G2").Value = "2016-10-17 11:50:07"
$dede.Range("H2").Value     = "2016-10-17 11:50:07"
$zhcn.Range("H2").Value     = "2016-10-17 11:49:56"

# --- Widen the Status / datetime columns on every sheet ----------------
$overview.Range("E1:F1").ColumnWidth = 16.382654825846366
$zhcn.Range("C1").ColumnWidth        = 16.382654825846366
$dede.Range("C1").ColumnWidth        = 16.382654825846366
